# "Generate Report for Handback" - refresh the handoff/handback timestamps
# recorded on the per-locale handback-status sheets.
$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("G2").Value = "2016-07-26 07:50:48"
$zhcn.Range("J2").Value = "2016-07-26 07:51:48"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("J2").Value = "2016-07-26 07:52:05"
